# Update gh-pages output data: increment "views" (column F) counters
# for a handful of rows in the "展览" (sheet1) and "全部类型" (sheet4)
# worksheets, mirroring a re-scrape of the source data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---------------------------------------------------
$wsExhibit = $wb.Worksheets.Item("展览")

$wsExhibit.Range("F3").Value  = 2931
$wsExhibit.Range("F4").Value  = 2931
$wsExhibit.Range("F5").Value  = 6406
$wsExhibit.Range("F6").Value  = 2512
$wsExhibit.Range("F8").Value  = 415
$wsExhibit.Range("F14").Value = 7316
$wsExhibit.Range("F40").Value = 731
$wsExhibit.Range("F42").Value = 11
$wsExhibit.Range("F43").Value = 197
$wsExhibit.Range("F46").Value = 5
$wsExhibit.Range("F49").Value = 50
$wsExhibit.Range("F50").Value = 41

# --- Sheet "全部类型" -------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")

$wsAll.Range("F4").Value  = 2931
$wsAll.Range("F9").Value  = 6406
$wsAll.Range("F19").Value = 7316
$wsAll.Range("F25").Value = 8877
$wsAll.Range("F42").Value = 731
$wsAll.Range("F44").Value = 197
$wsAll.Range("F49").Value = 50
$wsAll.Range("F50").Value = 41
